# Applies the "Penalty Reward System" forecast-shift edit:
#   - Forecast Comparison sheet: each week's Week_Start_Date (col B) is
#     shifted forward by one week, and MyForecast (col D) gets new values.
#   - Summary sheet: the derived metrics are refreshed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: columns B (Week_Start_Date) and D (MyForecast) ---
# Dates are entered with a leading apostrophe so Excel stores them as literal
# text (matching the workbook's existing inline-string date format) instead
# of auto-converting them to date serial numbers.

$rows = @(
    @{ Row = 2;  Date = "2025-01-12"; Forecast = 12 },
    @{ Row = 3;  Date = "2025-01-19"; Forecast = 11 },
    @{ Row = 4;  Date = "2025-01-26"; Forecast = 10 },
    @{ Row = 5;  Date = "2025-02-02"; Forecast = 10 },
    @{ Row = 6;  Date = "2025-02-09"; Forecast = 9  },
    @{ Row = 7;  Date = "2025-02-16"; Forecast = 8  },
    @{ Row = 8;  Date = "2025-02-23"; Forecast = 7  },
    @{ Row = 9;  Date = "2025-03-02"; Forecast = 6  },
    @{ Row = 10; Date = "2025-03-09"; Forecast = 6  },
    @{ Row = 11; Date = "2025-03-16"; Forecast = 14 },
    @{ Row = 12; Date = "2025-03-23"; Forecast = 14 },
    @{ Row = 13; Date = "2025-03-30"; Forecast = 16 },
    @{ Row = 14; Date = "2025-04-06"; Forecast = 14 },
    @{ Row = 15; Date = "2025-04-13"; Forecast = 13 },
    @{ Row = 16; Date = "2025-04-20"; Forecast = 14 },
    @{ Row = 17; Date = "2025-04-27"; Forecast = 13 }
)

foreach ($r in $rows) {
    $wsForecast.Range("B$($r.Row)").Value = "'" + $r.Date
    $wsForecast.Range("D$($r.Row)").Value = $r.Forecast
}

# --- Summary sheet updates ---
$wsSummary.Range("B2").Value  = "2023-01-22 to 2025-01-05"
$wsSummary.Range("B8").Value  = "2523 units"
$wsSummary.Range("B9").Value  = "'178"
$wsSummary.Range("B10").Value = "'73"
$wsSummary.Range("B11").Value = "'43"
$wsSummary.Range("B12").Value = "'16"
$wsSummary.Range("B13").Value = "'2025-03-30"
$wsSummary.Range("B14").Value = "'6"
$wsSummary.Range("B15").Value = "'2025-03-02"
